# Commit: "Update piTest to use DataSet"
#
# Semantic change: the "Time unit" column (K) on the Boswell_2012 sheet
# held the literal text "day[s]" for every data row. The data set was
# updated to spell the unit "day(s)" instead (bracket -> parens), and the
# workbook was left with the Boswell_2012 sheet active / selected at H2
# (previously Cedersund_2008 was the active tab).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Boswell_2012")

# "Time unit" values live in K2:K51 (K1 is the header "Time unit").
# Replace the old "day[s]" label with "day(s)" for every data row.
$ws.Range("K2:K51").Value2 = "day(s)"

# Make Boswell_2012 the active sheet and put the cursor on H2, matching
# the saved selection/view state after the edit.
$ws.Activate()
$ws.Range("H2").Select()
